# Auto-generated edit script for Regional.xlsx (Tables/Regional.xlsx)
# Commit message: "Had to update data with complete dataset back to 1980
# then reprocess."
#
# All touched cells store plain TEXT (this workbook was authored by
# openpyxl with inline strings) even though most values look numeric or
# date-like. A bare Range.Value = "4.24" would make Excel auto-detect a
# Number/Date, which (a) changes the stored value type and (b) swaps in a
# brand-new number-format style instead of reusing the sheet-wide style
# used by every data cell. So for every edit we:
#   1) enter the text with a leading apostrophe (the normal Excel way to
#      force literal text, exactly like typing it into the grid), then
#   2) formats-only paste from an untouched, same-style neighbour cell to
#      drop the quote-prefix style Excel just created and land back on
#      the original style index.

$wb = $excel.ActiveWorkbook
$wsMap = $wb.Worksheets.Item("Regional for Mapping")
$wsExhibit = $wb.Worksheets.Item("Regional Exhibit")

# --- Enter the new values as forced text ---------------------------------
$wsMap.Range("X2").Value = "'4.24"
$wsMap.Range("Y2").Value = "'2020-03-11"
$wsMap.Range("AB2").Value = "'3.6"
$wsMap.Range("AD2").Value = "'95"
$wsMap.Range("X3").Value = "'5.16"
$wsMap.Range("Y3").Value = "'2020-03-11"
$wsMap.Range("AB3").Value = "'3.7"
$wsMap.Range("AD3").Value = "'87"
$wsMap.Range("Y7").Value = "'2020-03-16"
$wsMap.Range("AD7").Value = "'127"
$wsMap.Range("X8").Value = "'5.05"
$wsMap.Range("Y8").Value = "'2020-03-17"
$wsMap.Range("AD8").Value = "'103"
$wsMap.Range("X9").Value = "'5.08"
$wsMap.Range("Y9").Value = "'2020-03-16"
$wsMap.Range("AD9").Value = "'96"
$wsMap.Range("AB15").Value = "'4.8"
$wsMap.Range("AD15").Value = "'53"
$wsMap.Range("AD16").Value = "'82"
$wsMap.Range("X17").Value = "'5.68"
$wsMap.Range("X19").Value = "'5.55"
$wsMap.Range("X31").Value = "'7.38"
$wsMap.Range("Y31").Value = "'2020-03-12"
$wsMap.Range("AD31").Value = "'98"
$wsMap.Range("X32").Value = "'8.04"
$wsMap.Range("Y32").Value = "'2020-03-12"
$wsMap.Range("AD32").Value = "'106"
$wsMap.Range("X33").Value = "'25.7"
$wsMap.Range("Y33").Value = "'2020-03-13"
$wsMap.Range("AD33").Value = "'94"
$wsMap.Range("V34").Value = "'47.6"
$wsMap.Range("W34").Value = "'2020-03-13"
$wsMap.Range("X34").Value = "'47.6"
$wsMap.Range("Y34").Value = "'2020-03-13"
$wsMap.Range("AD34").Value = "'89"
$wsMap.Range("X51").Value = "'5.29"
$wsMap.Range("Y51").Value = "'2020-03-18"
$wsMap.Range("AD51").Value = "'54"
$wsMap.Range("AD53").Value = "'79"
$wsMap.Range("X60").Value = "'2.08"
$wsMap.Range("AD65").Value = "'46"
$wsExhibit.Range("G37").Value = "'47.6"
$wsExhibit.Range("H37").Value = "'3/13/20"

# --- Restore the original (style 1) formatting on each edited cell -------
$wsMap.Range("A1").Copy()
$wsMap.Range("X2:Y2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AB2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X3:Y3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AB3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("Y7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X8:Y8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X9:Y9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AB15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X31:Y31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X32:Y32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X33:Y33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("V34:Y34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X51:Y51").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD51").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD53").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("X60").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsMap.Range("A1").Copy()
$wsMap.Range("AD65").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$wsExhibit.Range("A37").Copy()
$wsExhibit.Range("G37:H37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false

# --- Column width tweaks on "Regional Exhibit" (cols A and F) ------------
# The engine quantizes ColumnWidth to the nearest 1/6 of a character, so
# request the value whose quantized result lands on the authored width.
$wsExhibit.Columns.Item(1).ColumnWidth = 11.5   # -> stored width 12.25-ish (12 1/3)
$wsExhibit.Columns.Item(6).ColumnWidth = 15.5   # -> stored width 16.25-ish (16 1/3)
